$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit is a permutation of the data rows (2-13): each row's content
# (columns D, J, K, L, M, N, O, P, Q) moves to a different row, per the
# mapping derived from the diff. Capture the "before" snapshot first,
# then write the permuted values back so row->row moves don't clobber
# data that hasn't been read yet.

$mapping = @{
    2  = 11
    3  = 13
    4  = 5
    5  = 7
    6  = 9
    7  = 12
    8  = 6
    9  = 10
    10 = 3
    11 = 4
    12 = 8
    13 = 2
}

$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

$snapshot = @{}
foreach ($r in 2..13) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($r in 2..13) {
    $dest = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$dest").Value2 = $snapshot[$r][$c]
    }
}
